$wb = $excel.ActiveWorkbook

# --- 2.1.2 / 2.1.3: fix swapped title references ---------------------------
# Sheet "2.1.2" (Count of Teachers ...) incorrectly showed the "NUMBER OF
# SCHOOLS" title; it should read "NUMBER OF TEACHERS". Sheet "2.1.3" (Count
# of Enrolment ...) should show the "ENROLMENT" title that used to sit on
# "2.1.2".
$wsTeachers = $wb.Worksheets.Item("2.1.2")
$wsTeachers.Range("A1").Value = "2.1.1 NUMBER OF TEACHERS ACCORDING TO EDUCATION LEVEL BY DISTRICT"

$wsEnrolment = $wb.Worksheets.Item("2.1.3")
$wsEnrolment.Range("A1").Value = "2.1.3 ENROLMENT ACCORDING TO EDUCATION LEVEL BY DISTRICT"

# --- 2.2.2 / 2.2.3: replace "NA" placeholders with 0 counts -----------------
$ws222 = $wb.Worksheets.Item("2.2.2")
$ws222.Range("C4:D7").Value = 0

$ws223 = $wb.Worksheets.Item("2.2.3")
$ws223.Range("C4:D7").Value = 0

# --- View-state updates (selection / active cell per sheet) ----------------
# Apply these on the sheets that must stay inactive first, since selecting a
# range implicitly activates its worksheet.
$wb.Worksheets.Item("Data").Range("C3").Select() | Out-Null
$wb.Worksheets.Item("2.1.2").Range("I14").Select() | Out-Null
$wb.Worksheets.Item("2.1.5.2").Range("K7").Select() | Out-Null
$ws223.Range("C4:D7").Select() | Out-Null

# "2.2.2" ends up as the active sheet/tab, so select on it last.
$ws222.Range("G8").Select() | Out-Null

Write-Output "done"
